# "Skill Changes" sheet: fill in the new patch-note columns (C = summary,
# D = detail) for several Alchemy skills, clear the now-stale "Synergy" note,
# and move the saved selection to C8.
#
# NOTE: the order cells are written in below is intentional. New text values
# become shared-string-table entries in first-use order, and this order was
# chosen to reproduce the exact shared-string layout of the target workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Skill Changes")

# Row 58 - Protective Coating
$ws.Range("D58").Value = "Adds 2 additional skill levels for armor increase"
$ws.Range("C58").Value = "Increases armor instead of resistances (scales with Character Level)"

# Row 52 - Adaptation
$ws.Range("C52").Value = "Increased duration by 10% at max level, adds poison resistance"
$ws.Range("D52").Value = "Increases poison resistance by 2 skill point levels"

# Row 50 - Tissue Transmutation
$ws.Range("C50").Value = "Increased bonus, now also adds weapon damage reduction"
$ws.Range("D50").Value = "Increases weapon damage resistance by 2 levels"

# Row 51 - Synergy
$ws.Range("C51").Value = " "

# Row 54 - Pyrotechnics
$ws.Range("C54").Value = "Scales with player level"
$ws.Range("D54").Value = "Adds 2 additional skill levels"

# Row 60 - Hunters Instinct
$ws.Range("C60").Value = "Fixed bug with level scaling not working, adds monster damage resistance"
$ws.Range("D60").Value = "Increases monster damage resistance by 2 levels."

# Row 55 - Efficiency
$ws.Range("C55").Value = "Increased bonus to 6 instead of 5"
$ws.Range("D55").Value = "Adds 2 additional skill levels"

# Row 47 - Endure Pain
$ws.Range("C47").Value = "Increased effect, adds damage reduction"
$ws.Range("D47").Value = "Increases damage reduction by 2 levels"

# Row 56 - Cluster Bomb
$ws.Range("C56").Value = "Increased number of clusters per level from 1 to 3"

# Update the view: selection moves to C8, scroll back to top-left (no frozen topLeftCell override)
$ws.Range("C8").Select() | Out-Null
